$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sheet1"

# Header row (row 1) values
$ws.Cells.Item(1,1).Value = "Date"
$ws.Cells.Item(1,2).Value = "Model Name"
$ws.Cells.Item(1,3).Value = "Exact Precision (Micro Avg)"
$ws.Cells.Item(1,4).Value = "Exact Recall (Micro Avg)"
$ws.Cells.Item(1,5).Value = "Exact F1 Score (Micro Avg)"
$ws.Cells.Item(1,6).Value = "Exact Precision (Macro Avg)"
$ws.Cells.Item(1,7).Value = "Exact Recall (Macro Avg)"
$ws.Cells.Item(1,8).Value = "Exact F1 Score (Macro Avg)"
$ws.Cells.Item(1,9).Value = "Exact Precision (Weighted Avg)"
$ws.Cells.Item(1,10).Value = "Exact Recall (Weighted Avg)"
$ws.Cells.Item(1,11).Value = "Exact F1 Score (Weighted Avg)"
$ws.Cells.Item(1,12).Value = "Partial Precision"
$ws.Cells.Item(1,13).Value = "Partial Recall"
$ws.Cells.Item(1,14).Value = "Partial F1 Score"
$ws.Cells.Item(1,15).Value = "Partial TP"
$ws.Cells.Item(1,16).Value = "Partial FP"
$ws.Cells.Item(1,17).Value = "Partial FN"
$ws.Cells.Item(1,18).Value = "Support"
$ws.Cells.Item(1,19).Value = "Accuracy"
$ws.Cells.Item(1,20).Value = "Result Link"
$ws.Cells.Item(1,21).Value = "Stats Link"
$ws.Cells.Item(1,22).Value = "No of GPU Used"
$ws.Cells.Item(1,23).Value = "Power Consumption"
$ws.Cells.Item(1,24).Value = "Unnamed: 23"

# Apply header style: bold font, thin box border, center/top alignment
$hdr = $ws.Range("A1:X1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

# Row 2 data
$ws.Cells.Item(2,1).Value = "'09/11/2025"
$ws.Cells.Item(2,2).Value = "Qwen2.5-32B-Instruct"
$ws.Cells.Item(2,3).Value = 0.4642857142857143
$ws.Cells.Item(2,4).Value = 0.3063973063973064
$ws.Cells.Item(2,5).Value = 0.3691683569979716
$ws.Cells.Item(2,6).Value = 0.2226114348853367
$ws.Cells.Item(2,7).Value = 0.1434872517929153
$ws.Cells.Item(2,8).Value = 0.169306566378869
$ws.Cells.Item(2,9).Value = 0.482058700890429
$ws.Cells.Item(2,10).Value = 0.3063973063973064
$ws.Cells.Item(2,11).Value = 0.3658240415704264
$ws.Cells.Item(2,12).Value = 0.5721649484536082
$ws.Cells.Item(2,13).Value = 0.375
$ws.Cells.Item(2,14).Value = 0.4530612244897959
$ws.Cells.Item(2,15).Value = 111
$ws.Cells.Item(2,16).Value = 83
$ws.Cells.Item(2,17).Value = 185
$ws.Cells.Item(2,18).Value = 297
$ws.Cells.Item(2,19).Value = 0.9563810665068904
$ws.Cells.Item(2,20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_3_shot.txt"
$ws.Cells.Item(2,21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_3_shot.txt"
$ws.Cells.Item(2,22).Value = "4 MLGPU"
$ws.Cells.Item(2,23).Value = "0.098 kWh"
$ws.Cells.Item(2,24).Value = 3168

# Row 3 data
$ws.Cells.Item(3,1).Value = "'09/12/2025"
$ws.Cells.Item(3,2).Value = "Qwen2.5-32B-Instruct"
$ws.Cells.Item(3,3).Value = 0.4717741935483871
$ws.Cells.Item(3,4).Value = 0.3939393939393939
$ws.Cells.Item(3,5).Value = 0.4293577981651376
$ws.Cells.Item(3,6).Value = 0.5566959308807135
$ws.Cells.Item(3,7).Value = 0.3903227756597969
$ws.Cells.Item(3,8).Value = 0.4449845271973715
$ws.Cells.Item(3,9).Value = 0.5823406866885128
$ws.Cells.Item(3,10).Value = 0.3939393939393939
$ws.Cells.Item(3,11).Value = 0.455253739394065
$ws.Cells.Item(3,12).Value = 0.5161290322580645
$ws.Cells.Item(3,13).Value = 0.4324324324324325
$ws.Cells.Item(3,14).Value = 0.4705882352941176
$ws.Cells.Item(3,15).Value = 128
$ws.Cells.Item(3,16).Value = 120
$ws.Cells.Item(3,17).Value = 168
$ws.Cells.Item(3,18).Value = 297
$ws.Cells.Item(3,19).Value = 0.9458358298382264
$ws.Cells.Item(3,20).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/ner_evaluation_results_Qwen2.5-32B-Instruct_3_shot.txt"
$ws.Cells.Item(3,21).Value = "/home/s27mhusa_hpc/Master-Thesis/Evaluation_Results/Final_TestFiles_3rdSeptember_FewShotTest_Broad/Stats/ner_evaluation_stats_Qwen2.5-32B-Instruct_3_shot.txt"
$ws.Cells.Item(3,22).Value = "4 MLGPU"
$ws.Cells.Item(3,23).Value = "0.061 kWh"
$ws.Cells.Item(3,24).Value = "'"
